$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Find the last populated row in column A and append the new book right
# after it (mirrors the diff: a single new row 60 appended to the table).
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Copy the formatting (incl. the date number format) down from the row
# above so the new row reuses the existing style ids instead of minting
# new ones.
$ws.Range("A" + $lastRow + ":G" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":G" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = "The Little Book of Common Sense Investing"
$ws.Cells.Item($newRow, 2).Value = "John C. Bogle"
$ws.Cells.Item($newRow, 3).Value = 43942
$ws.Cells.Item($newRow, 4).Value = 43942
$ws.Cells.Item($newRow, 5).Value = "investing;stock market;index fund;passive investing"
$ws.Cells.Item($newRow, 6).Value = "Audio"
$ws.Cells.Item($newRow, 7).Value = "5 Hours 7 Mins"

$ws.Cells.Item($newRow + 1, 1).Select()
